# Auto-generated Excel COM-interop script
# Applies market-data refresh values (scheduled runner update) to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 308280.2
$ws.Range("I15").Value = 308280.2
$ws.Range("K15").Value = 924840.6000000001
$ws.Range("M15").Value = -924671.6000000001
$ws.Range("H19").Value = 1951.8286
$ws.Range("I19").Value = 3700.8667
$ws.Range("J19").Value = 640.05
$ws.Range("K19").Value = 3700.8667
$ws.Range("L19").Value = 640.05
$ws.Range("M19").Value = -3525.8667
$ws.Range("N19").Value = -990.05
$ws.Range("H32").Value = 3319444
$ws.Range("I32").Value = 461.5
$ws.Range("J32").Value = 5361895
$ws.Range("K32").Value = 461.5
$ws.Range("L32").Value = 5361895
$ws.Range("M32").Value = -135.5
$ws.Range("N32").Value = -5362547
$ws.Range("H68").Value = 41333.332
$ws.Range("J68").Value = 41333.332
$ws.Range("L68").Value = 41333.332
$ws.Range("N68").Value = -42831.332
$ws.Range("H71").Value = 41333.332
$ws.Range("J71").Value = 41333.332
$ws.Range("L71").Value = 123999.996
$ws.Range("N71").Value = -131487.996
$ws.Range("H116").Value = 4688.4
$ws.Range("I116").Value = 5485.3335
$ws.Range("J116").Value = 3493
$ws.Range("K116").Value = 5485.3335
$ws.Range("L116").Value = 3493
$ws.Range("M116").Value = -2043.3335
$ws.Range("N116").Value = -10377
$ws.Range("H137").Value = 3763.875
$ws.Range("I137").Value = 5742.8335
$ws.Range("J137").Value = 1784.9166
$ws.Range("K137").Value = 17228.5005
$ws.Range("L137").Value = 5354.7498
$ws.Range("M137").Value = -14678.5005
$ws.Range("N137").Value = -10454.7498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1551.875
$ws.Range("I63").Value = 1616
$ws.Range("J63").Value = 1103
$ws.Range("K63").Value = 1616
$ws.Range("L63").Value = 1103
$ws.Range("M63").Value = -930
$ws.Range("N63").Value = -2475
$ws.Range("H66").Value = 1551.875
$ws.Range("I66").Value = 1616
$ws.Range("J66").Value = 1103
$ws.Range("K66").Value = 8080
$ws.Range("L66").Value = 5515
$ws.Range("M66").Value = -4648
$ws.Range("N66").Value = -12379
$ws.Range("H132").Value = 21587.389
$ws.Range("I132").Value = 32196.059
$ws.Range("J132").Value = 3552.65
$ws.Range("K132").Value = 96588.177
$ws.Range("L132").Value = 10657.95
$ws.Range("M132").Value = -94058.177
$ws.Range("N132").Value = -15717.95

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2732.6123
$ws.Range("I31").Value = 1869.9706
$ws.Range("J31").Value = 4687.933
$ws.Range("K31").Value = 1869.9706
$ws.Range("L31").Value = 4687.933
$ws.Range("M31").Value = -1574.9706
$ws.Range("N31").Value = -5277.933
$ws.Range("H34").Value = 2732.6123
$ws.Range("I34").Value = 1869.9706
$ws.Range("J34").Value = 4687.933
$ws.Range("K34").Value = 1869.9706
$ws.Range("L34").Value = 4687.933
$ws.Range("M34").Value = -1667.9706
$ws.Range("N34").Value = -5091.933
$ws.Range("H99").Value = 38834.11
$ws.Range("I99").Value = 60390.176
$ws.Range("J99").Value = 2188.8
$ws.Range("K99").Value = 60390.176
$ws.Range("L99").Value = 2188.8
$ws.Range("M99").Value = -58892.176
$ws.Range("N99").Value = -5184.8
$ws.Range("H126").Value = 38834.11
$ws.Range("I126").Value = 60390.176
$ws.Range("J126").Value = 2188.8
$ws.Range("K126").Value = 181170.528
$ws.Range("L126").Value = 6566.400000000001
$ws.Range("M126").Value = -178700.528
$ws.Range("N126").Value = -11506.4
$ws.Range("H132").Value = 2248.44
$ws.Range("I132").Value = 1387.5333
$ws.Range("J132").Value = 3539.8
$ws.Range("K132").Value = 4162.5999
$ws.Range("L132").Value = 10619.4
$ws.Range("M132").Value = -1632.5999
$ws.Range("N132").Value = -15679.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 992.3721
$ws.Range("I5").Value = 394.48
$ws.Range("J5").Value = 1822.7778
$ws.Range("K5").Value = 1183.44
$ws.Range("L5").Value = 5468.3334
$ws.Range("M5").Value = -1071.44
$ws.Range("N5").Value = -5692.3334
$ws.Range("H33").Value = 523.4
$ws.Range("I33").Value = 420
$ws.Range("J33").Value = 575.1
$ws.Range("K33").Value = 2520
$ws.Range("L33").Value = 3450.6
$ws.Range("M33").Value = -2237
$ws.Range("N33").Value = -4016.6
$ws.Range("H86").Value = 679.93335
$ws.Range("I86").Value = 678.5714
$ws.Range("J86").Value = 699
$ws.Range("K86").Value = 2035.7142
$ws.Range("L86").Value = 2097
$ws.Range("M86").Value = -849.7142000000001
$ws.Range("N86").Value = -4469
$ws.Range("H89").Value = 679.93335
$ws.Range("I89").Value = 678.5714
$ws.Range("J89").Value = 699
$ws.Range("K89").Value = 6107.1426
$ws.Range("L89").Value = 6291
$ws.Range("M89").Value = -179.1426000000001
$ws.Range("N89").Value = -18147
$ws.Range("H131").Value = 1184.0944
$ws.Range("J131").Value = 1356.5227
$ws.Range("L131").Value = 4069.5681
$ws.Range("N131").Value = -14149.5681
$ws.Range("H135").Value = 992.3721
$ws.Range("I135").Value = 394.48
$ws.Range("J135").Value = 1822.7778
$ws.Range("K135").Value = 3550.32
$ws.Range("L135").Value = 16405.0002
$ws.Range("M135").Value = -1015.32
$ws.Range("N135").Value = -21475.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4120.743
$ws.Range("I70").Value = 3931.1592
$ws.Range("J70").Value = 4398.8
$ws.Range("K70").Value = 3931.1592
$ws.Range("L70").Value = 4398.8
$ws.Range("M70").Value = -3661.1592
$ws.Range("N70").Value = -4938.8
$ws.Range("H73").Value = 4120.743
$ws.Range("I73").Value = 3931.1592
$ws.Range("J73").Value = 4398.8
$ws.Range("K73").Value = 3931.1592
$ws.Range("L73").Value = 4398.8
$ws.Range("M73").Value = -2995.1592
$ws.Range("N73").Value = -6270.8
$ws.Range("H80").Value = 3807.6
$ws.Range("I80").Value = 4302.212
$ws.Range("J80").Value = 2847.4707
$ws.Range("K80").Value = 4302.212
$ws.Range("L80").Value = 2847.4707
$ws.Range("M80").Value = -3304.212
$ws.Range("N80").Value = -4843.4707
$ws.Range("H83").Value = 3807.6
$ws.Range("I83").Value = 4302.212
$ws.Range("J83").Value = 2847.4707
$ws.Range("K83").Value = 21511.06
$ws.Range("L83").Value = 14237.3535
$ws.Range("M83").Value = -16519.06
$ws.Range("N83").Value = -24221.3535
$ws.Range("H107").Value = 6225.4707
$ws.Range("I107").Value = 8654.583000000001
$ws.Range("J107").Value = 395.6
$ws.Range("K107").Value = 8654.583000000001
$ws.Range("L107").Value = 395.6
$ws.Range("M107").Value = -6734.583000000001
$ws.Range("N107").Value = -4235.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 649.55554
$ws.Range("I22").Value = 448.5
$ws.Range("J22").Value = 750.0833
$ws.Range("K22").Value = 448.5
$ws.Range("L22").Value = 750.0833
$ws.Range("M22").Value = -153.5
$ws.Range("N22").Value = -1340.0833
$ws.Range("H27").Value = 649.55554
$ws.Range("I27").Value = 448.5
$ws.Range("J27").Value = 750.0833
$ws.Range("K27").Value = 448.5
$ws.Range("L27").Value = 750.0833
$ws.Range("M27").Value = -341.5
$ws.Range("N27").Value = -964.0833
$ws.Range("H61").Value = 1600.8889
$ws.Range("I61").Value = 1384.6666
$ws.Range("J61").Value = 2033.3334
$ws.Range("K61").Value = 1384.6666
$ws.Range("L61").Value = 2033.3334
$ws.Range("M61").Value = -1182.6666
$ws.Range("N61").Value = -2437.3334
$ws.Range("H113").Value = 1600.8889
$ws.Range("I113").Value = 1384.6666
$ws.Range("J113").Value = 2033.3334
$ws.Range("K113").Value = 1384.6666
$ws.Range("L113").Value = 2033.3334
$ws.Range("M113").Value = 785.3334
$ws.Range("N113").Value = -6373.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14444668
$ws.Range("I136").Value = 21299226
$ws.Range("J136").Value = 437527.4
$ws.Range("K136").Value = 63897678
$ws.Range("L136").Value = 1312582.2
$ws.Range("M136").Value = -63895128
$ws.Range("N136").Value = -1317682.2
